$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.36247181892395
$ws.Range("B1").Value = 2.681127071380615
$ws.Range("C1").Value = 2.723663568496704
$ws.Range("D1").Value = 3.446640253067017
$ws.Range("E1").Value = 1.861255764961243
